$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 347
$ws.Range("D2").Value = 53
$ws.Range("B5").Value = 0.8675
$ws.Range("D5").Value = 0.1325
